$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.103.27"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "2.497.31"
$ws.Range("E3").Value = "  -1.01%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'319.76"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("D6").Value = "'105.85"
$ws.Range("E6").Value = "  -3.00%  "
$ws.Range("E7").Value = "  -1.18%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -4.22%  "
$ws.Range("D10").Value = "'38.82"
$ws.Range("E10").Value = "  -4.01%  "
$ws.Range("D11").Value = "'20.03"
$ws.Range("E11").Value = "  -0.54%  "
$ws.Range("D12").Value = "'0.0802"
$ws.Range("E12").Value = "  -2.10%  "
$ws.Range("E14").Value = "  -2.26%  "
$ws.Range("D15").Value = "2.889.43"
$ws.Range("E15").Value = "  -0.95%  "
$ws.Range("D16").Value = "2.501.72"
$ws.Range("E16").Value = "  -1.04%  "
$ws.Range("D17").Value = "'0.831"
$ws.Range("E17").Value = "  -3.42%  "
$ws.Range("D18").Value = "47.950.93"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("E19").Value = "  -1.52%  "
$ws.Range("E20").Value = "  +9.07%  "
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").Value = "0.0₃0931"
$ws.Range("E22").Value = "  -1.52%  "
$ws.Range("D23").Value = "'71.09"
$ws.Range("E23").Value = "  -1.97%  "
$ws.Range("D24").Value = "'271.24"
$ws.Range("E24").Value = "  +1.05%  "
$ws.Range("D25").Value = "'2.51"
$ws.Range("E25").Value = "  -2.39%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'25.81"
$ws.Range("E27").Value = "  -1.56%  "
$ws.Range("D28").Value = "'2.29"
$ws.Range("E28").Value = "  -0.72%  "
$ws.Range("E29").Value = "  -4.27%  "
$ws.Range("D30").Value = "'0.140"
$ws.Range("E30").Value = "  -3.72%  "
$ws.Range("D31").Value = "'34.83"
$ws.Range("E31").Value = "  -0.41%  "
$ws.Range("D32").Value = "'49.34"
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("D34").Value = "'19.07"
$ws.Range("E34").Value = "  -4.76%  "
$ws.Range("D36").Value = "'0.0772"
$ws.Range("E36").Value = "  -2.49%  "
$ws.Range("E37").Value = "  -2.60%  "
$ws.Range("E38").Value = "  -3.04%  "
$ws.Range("E39").Value = "  -4.56%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").Value = "'22.47"
$ws.Range("E40").Value = "  +1.37%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "'121.95"
$ws.Range("E41").Value = "  +2.60%  "
$ws.Range("E42").Value = "  -2.21%  "
$ws.Range("E43").Value = "  +1.20%  "
$ws.Range("D44").Value = "'0.0302"
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("D45").Value = "1.998.46"
$ws.Range("E45").Value = "  -0.18%  "
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("E47").Value = "  -1.24%  "
$ws.Range("E48").Value = "  -1.22%  "
$ws.Range("D49").Value = "'8.91"
$ws.Range("E49").Value = "  -1.93%  "
$ws.Range("E50").Value = "  -1.94%  "
$ws.Range("D51").Value = "'79.02"
$ws.Range("E51").Value = "  -1.81%  "
